$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6's Category cell was a typo ("Strins"); point it at the already-correct
# "Strings" text used elsewhere in the sheet (row 4) instead of the typo.
$ws.Range("D6").Value = "Strings"

# New solution added: #682 "BaseballGame", solved with a Stack, in Python.
$ws.Range("B7").Value = 682
$ws.Range("C7").Value = "BaseballGame "
$ws.Range("D7").Value = "Stack "
$ws.Range("E7").Value = "Python"

# Rows 8-11 had picked up a slightly different (but visually identical) "no
# fill" style than the rest of the filled-in rows. Re-apply the same "no
# fill" interior so they settle back onto the common style used by rows 4-7.
$ws.Range("B8:B11").Interior.Pattern = -4142

# After typing the new row, Excel's selection sits on the next empty row.
$ws.Range("B8").Select()
